$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- At-bat 1 (rows 10-12) ---
$ws.Range("F10").Value = "CB"
$ws.Range("G10").Value = "Take"
$ws.Range("H10").Value = "Ball"
$ws.Range("M10").Value = "77.22 MPH"

$ws.Range("F11").Value = "FB"
$ws.Range("G11").Value = "Take"
$ws.Range("H11").Value = "Ball"

$ws.Range("F12").Value = "FB"
$ws.Range("G12").Value = "Swing"
$ws.Range("H12").Value = "In Play"
$ws.Range("M12").Value = "3.16°"

$ws.Range("J17").Value = "CH,CB,FB"

# --- At-bat 2 (rows 19-21) ---
$ws.Range("F19").Value = "FB"
$ws.Range("G19").Value = "Take"
$ws.Range("H19").Value = "Ball"
$ws.Range("M19").Value = "88.41 MPH"

$ws.Range("F20").Value = "FB"
$ws.Range("G20").Value = "Swing"
$ws.Range("H20").Value = "In Play"

$ws.Range("M21").Value = "44.38°"

$ws.Range("J26").Value = "CH,CB,FB"

# --- At-bat 3 (rows 28-31) ---
$ws.Range("F28").Value = "FB"
$ws.Range("G28").Value = "Take"
$ws.Range("H28").Value = "Ball"
$ws.Range("M28").Value = "98.92 MPH"

$ws.Range("F29").Value = "FB"
$ws.Range("G29").Value = "Take"
$ws.Range("H29").Value = "Strike"

$ws.Range("F30").Value = "CB"
$ws.Range("G30").Value = "Take"
$ws.Range("H30").Value = "Ball"
$ws.Range("M30").Value = "37.45°"

$ws.Range("F31").Value = "FB"
$ws.Range("G31").Value = "Swing"
$ws.Range("H31").Value = "In Play"

$ws.Range("J35").Value = "CH,CB,FB,SL"

# --- At-bat 4 (rows 37-40) ---
$ws.Range("F37").Value = "CH"
$ws.Range("G37").Value = "Take"
$ws.Range("H37").Value = "Strike"
$ws.Range("M37").Value = "75.55 MPH"

$ws.Range("F38").Value = "CB"
$ws.Range("G38").Value = "Take"
$ws.Range("H38").Value = "Ball"

$ws.Range("F39").Value = "CH"
$ws.Range("G39").Value = "Take"
$ws.Range("H39").Value = "Ball"
$ws.Range("M39").Value = "15.16°"

$ws.Range("F40").Value = "CH"
$ws.Range("G40").Value = "Swing"
$ws.Range("H40").Value = "In Play"

$ws.Range("J44").Value = "CH,CB,FB,SL"
